$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$wsMeta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
$wsMeta.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# --- Include ValueSet #0 sheet updates ---
$wsInclude = $wb.Worksheets.Item("Include ValueSet #0")

# Append version timestamp to the two ValueSet URLs
$wsInclude.Range("A2").Value = "https://mos.esante.gouv.fr/NOS/JDV_J283-PrestationsIndirects_SERAFIN/FHIR/JDV-J283-PrestationsIndirects-SERAFIN|20241025120000"
$wsInclude.Range("A3").Value = "https://mos.esante.gouv.fr/NOS/JDV_J284-PrestationsDirects_SERAFIN/FHIR/JDV-J284-PrestationsDirects-SERAFIN|20241025120000"

$wb.Save()
